$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates are Excel serial numbers). Column A keeps
# the same cell style ("s=2": bold, bordered, centered, custom date
# numFmt) as the preceding rows, so copy formatting from the row above
# before writing the new value.
$newRows = @(
    @{ Row = 252; A = 44326; B = 2; C = 43; D = 130.3701907043022 },
    @{ Row = 253; A = 44327; B = 1; C = 40; D = 121.2745960040021 },
    @{ Row = 254; A = 44328; B = 1; C = 41; D = 124.3064609041021 },
    @{ Row = 255; A = 44329; B = 1; C = 26; D = 78.82848740260134 }
)

foreach ($r in $newRows) {
    $srcA = $ws.Cells.Item($r.Row - 1, 1)
    $dstA = $ws.Cells.Item($r.Row, 1)
    $srcA.Copy($dstA)
    $dstA.Value = $r.A

    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
